# The commit swaps the XML content that lives in the two theme parts of the
# package (ppt/theme/theme1.xml <-> ppt/theme/theme2.xml): what used to be
# the "Office Theme" colour scheme moves into the part used by the
# slide master (so the deck's live/visible theme becomes "Office Theme"
# coloured) and what used to be the "Integral" colour scheme moves into the
# other theme part. The two themes only differ in their <a:clrScheme> (12
# colours) - the font scheme and format scheme are identical between them -
# so the edit is expressed here as re-pointing the presentation's theme
# colour scheme at the "Office Theme" palette.

function HexToOle([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), in
# the exact order exposed by ThemeColorScheme.Item(1..12).
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$scheme = $s.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $scheme.Item($i).RGB = HexToOle($officeThemeColors[$i - 1])
}
